# CU-08 Ver proveedor - corrections per "correcciones lunes de 23 de septiembre"
#
# Strategy: each touched paragraph is replaced in full via Range.InsertXML with a
# fragment built from the original paragraph's own <w:p> attributes (to preserve
# w14:paraId/rsid/pPr etc., none of which the diff changes) plus the new run
# layout demanded by the diff.

$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Find-ParagraphIndexByText($doc, $marker) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        if ($t -like "*$marker*") {
            return $i
        }
    }
    return -1
}

function Replace-Paragraph($doc, $marker, $innerXml) {
    $idx = Find-ParagraphIndexByText $doc $marker
    if ($idx -eq -1) {
        Write-Host "WARNING: paragraph not found for marker: $marker"
        return
    }
    $para = $doc.Paragraphs($idx)
    $range = $para.Range
    $range.InsertXML($innerXml)
}

# --- Change 1 : Precondiciones / PRE-01 -----------------------------------
$p17 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="5DFC1C47" w14:textId="51CF949F" w:rsidR="00571977" w:rsidRDefault="00571977" w:rsidP="00244463">' `
  + '<w:pPr><w:jc w:val="both"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">PRE-01 </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">Debe existir por lo menos </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">un </w:t></w:r>' `
  + '<w:r><w:t>PROVEEDOR registrado en el sistema</w:t></w:r>' `
  + '</w:p>'
Replace-Paragraph $d "PRE-01" $p17

# --- Change 2 : Flujo normal / paso 1 (recupera de la base de datos...) ---
$p20 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="3536A61C" w14:textId="0C8F4928" w:rsidR="007D75ED" w:rsidRDefault="00571977" w:rsidP="00571977">' `
  + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">El sistema </w:t></w:r>' `
  + '<w:r><w:t>recupera de</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> la base de datos toda la información del PROVEEDOR</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> y los </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:t>PRODUCTO</w:t></w:r>' `
  + '<w:r><w:t>s</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:t>que vende</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">. </w:t></w:r>' `
  + '<w:r><w:t>(EX-01)</w:t></w:r>' `
  + '</w:p>'
Replace-Paragraph $d "consulta en la base de datos" $p20

# --- Change 3 : Flujo normal / texto "Luego muestra..." campos ------------
$p21 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="57A0C0A8" w14:textId="17FBE2C0" w:rsidR="00383506" w:rsidRDefault="00FB226B" w:rsidP="007D75ED">' `
  + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:jc w:val="both"/></w:pPr>' `
  + '<w:r><w:t>Luego m</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">uestra la ventana </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:t>DetalleProveedorView</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:t>con los campos correo, nombre, RFC</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">, </w:t></w:r>' `
  + '<w:r><w:t>teléfono</w:t></w:r>' `
  + '<w:r><w:t>, estado.</w:t></w:r>' `
  + '</w:p>'
Replace-Paragraph $d "con los campos correo, nombre, RFC" $p21

# --- Change 4 : Flujo normal / "Debajo, muestra una tabla..." -------------
$p22 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="580CE85B" w14:textId="198D38C2" w:rsidR="00383506" w:rsidRDefault="00383506" w:rsidP="00383506">' `
  + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:jc w:val="both"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">Debajo, muestra una tabla con los </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:t>PRODUCTO</w:t></w:r>' `
  + '<w:r><w:t>s</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> vendidos por el PROVEEDOR.</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> La tabla tiene los campos: código, nombre y descripción.</w:t></w:r>' `
  + '</w:p>'
Replace-Paragraph $d "Debajo, muestra una tabla" $p22

# --- Change 5 : Excepciones / EX-01 ----------------------------------------
$p32 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="56C940EF" w14:textId="77777777" w:rsidR="00383506" w:rsidRDefault="00383506" w:rsidP="00383506">' `
  + '<w:pPr><w:jc w:val="both"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">EX-01 </w:t></w:r>' `
  + '<w:r><w:t>No hay conexión a la red</w:t></w:r>' `
  + '</w:p>'
Replace-Paragraph $d "EX-01 No hay" $p32

Write-Host "All changes applied."
